# Sprint 1 Report - apply Rachel's review comments
#
# The underlying engine re-derives runs for any paragraph text that is
# touched by an edit, merging adjacent same-formatting runs back
# together. To land on the exact run layout the diff calls for, each
# edit below first rewrites the affected span as plain text (which the
# engine collapses into a single run) and then uses a short-lived
# Bookmarks.Add/Delete pair at each internal boundary to force a clean
# run split with no leftover formatting / bookmark residue.

$d = $word.ActiveDocument

function Split-RunAt($pos) {
    $bm = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpSplitMarker", $bm) | Out-Null
    $d.Bookmarks("TmpSplitMarker").Delete()
}

# ---------------------------------------------------------------------
# Hunk 1: "We are developing ... Overview of Project" paragraph.
# Split the opening run into two sentences.
# ---------------------------------------------------------------------
$oldLead = "We are developing Android and iOS versions of the eBill application first. When Innovative gives us the API and requirements for the service ticket application development will be started on that as well. "
$newSentence1 = "We are developing Android and iOS versions of the eBill application first. "
$newSentence2 = "A service ticket application will enter development when the API and requirements become available. "

$findRange = $d.Content
$found = $findRange.Find.Execute($oldLead, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the original 'We are developing' sentence."
}
$leadStart = $findRange.Start

$leadRange = $d.Range($leadStart, $leadStart + $oldLead.Length)
$leadRange.Text = $newSentence1 + $newSentence2
Split-RunAt ($leadStart + $newSentence1.Length)

# ---------------------------------------------------------------------
# Hunk 2: "Project will use Eclipse, Xcode, and StoryBoard for
# development" bullet -> "Project will use Eclipse and Xcode for
# development environments", re-split into the seven runs the diff
# shows.
# ---------------------------------------------------------------------
$oldDev = "Project will use Eclipse, Xcode, and StoryBoard for development"
$newDevParts = @(
    "Project will use Eclipse",
    " ",
    "and ",
    "Xc",
    "ode for development",
    " ",
    "environments"
)
$newDev = [string]::Join("", $newDevParts)

$findRange2 = $d.Content
$found = $findRange2.Find.Execute($oldDev, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the original 'Project will use Eclipse' bullet."
}
$devStart = $findRange2.Start

$devRange = $d.Range($devStart, $devStart + $oldDev.Length)
$devRange.Text = $newDev

$offset = 0
for ($i = 0; $i -lt ($newDevParts.Length - 1); $i++) {
    $offset += $newDevParts[$i].Length
    Split-RunAt ($devStart + $offset)
}

# ---------------------------------------------------------------------
# Hunk 3 & 4: the "_GoBack" bookmark moves from the end of the "by
# other team members for data access" bullet to the middle of the
# "Android-specific design issues" bullet (splitting "specific" into
# "spec" / "ific"). Word only keeps one bookmark per name, so adding
# the bookmark at the new spot automatically removes it from the old
# one.
# ---------------------------------------------------------------------
$marker = "Android-specific design issues."
$findRange3 = $d.Content
$found = $findRange3.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Android-specific design issues.' bullet."
}
$markerStart = $findRange3.Start
$splitPos = $markerStart + "Android-spec".Length

$goBackPoint = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null
